# Auto-generated edit script applying Odin_Profits.xlsx value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Cells.Item(64, 8).Value = 33337298  # H64
$ws.Cells.Item(64, 9).Value = 55557616  # I64
$ws.Cells.Item(64, 10).Value = 6825  # J64
$ws.Cells.Item(64, 11).Value = 55557616  # K64
$ws.Cells.Item(64, 12).Value = 6825  # L64
$ws.Cells.Item(64, 13).Value = -55557368  # M64
$ws.Cells.Item(64, 14).Value = -7321  # N64

# Row 67
$ws.Cells.Item(67, 8).Value = 33337298  # H67
$ws.Cells.Item(67, 9).Value = 55557616  # I67
$ws.Cells.Item(67, 10).Value = 6825  # J67
$ws.Cells.Item(67, 11).Value = 55557616  # K67
$ws.Cells.Item(67, 12).Value = 6825  # L67
$ws.Cells.Item(67, 13).Value = -55556758  # M67
$ws.Cells.Item(67, 14).Value = -8541  # N67

# Row 112
$ws.Cells.Item(112, 8).Value = 2107.9814  # H112
$ws.Cells.Item(112, 9).Value = 949.6  # I112
$ws.Cells.Item(112, 10).Value = 2226.1836  # J112
$ws.Cells.Item(112, 11).Value = 2848.8  # K112
$ws.Cells.Item(112, 12).Value = 6678.550799999999  # L112
$ws.Cells.Item(112, 13).Value = -1740.8  # M112
$ws.Cells.Item(112, 14).Value = -8894.550799999999  # N112

# Row 135
$ws.Cells.Item(135, 8).Value = 10452.444  # H135
$ws.Cells.Item(135, 9).Value = 6274.6665  # I135
$ws.Cells.Item(135, 10).Value = 12541.333  # J135
$ws.Cells.Item(135, 11).Value = 56471.9985  # K135
$ws.Cells.Item(135, 12).Value = 112871.997  # L135
$ws.Cells.Item(135, 13).Value = -53936.9985  # M135
$ws.Cells.Item(135, 14).Value = -117941.997  # N135

# Row 141
$ws.Cells.Item(141, 8).Value = 2510.8  # H141
$ws.Cells.Item(141, 9).Value = 1889.125  # I141
$ws.Cells.Item(141, 10).Value = 4997.5  # J141
$ws.Cells.Item(141, 11).Value = 5667.375  # K141
$ws.Cells.Item(141, 12).Value = 14992.5  # L141
$ws.Cells.Item(141, 13).Value = -487.375  # M141
$ws.Cells.Item(141, 14).Value = -25352.5  # N141


$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 1635.7142  # H45
$ws.Cells.Item(45, 9).Value = 1590  # I45
$ws.Cells.Item(45, 10).Value = 1750  # J45
$ws.Cells.Item(45, 11).Value = 1590  # K45
$ws.Cells.Item(45, 12).Value = 1750  # L45
$ws.Cells.Item(45, 13).Value = -1213  # M45
$ws.Cells.Item(45, 14).Value = -2504  # N45

# Row 63
$ws.Cells.Item(63, 8).Value = 5426.591  # H63
$ws.Cells.Item(63, 9).Value = 4608  # I63
$ws.Cells.Item(63, 10).Value = 5733.5625  # J63
$ws.Cells.Item(63, 11).Value = 4608  # K63
$ws.Cells.Item(63, 12).Value = 5733.5625  # L63
$ws.Cells.Item(63, 13).Value = -3922  # M63
$ws.Cells.Item(63, 14).Value = -7105.5625  # N63

# Row 66
$ws.Cells.Item(66, 8).Value = 5426.591  # H66
$ws.Cells.Item(66, 9).Value = 4608  # I66
$ws.Cells.Item(66, 10).Value = 5733.5625  # J66
$ws.Cells.Item(66, 11).Value = 23040  # K66
$ws.Cells.Item(66, 12).Value = 28667.8125  # L66
$ws.Cells.Item(66, 13).Value = -19608  # M66
$ws.Cells.Item(66, 14).Value = -35531.8125  # N66

# Row 74
$ws.Cells.Item(74, 8).Value = 5081.409  # H74
$ws.Cells.Item(74, 9).Value = 6363.4546  # I74
$ws.Cells.Item(74, 10).Value = 3799.3635  # J74
$ws.Cells.Item(74, 11).Value = 6363.4546  # K74
$ws.Cells.Item(74, 12).Value = 3799.3635  # L74
$ws.Cells.Item(74, 13).Value = -5489.4546  # M74
$ws.Cells.Item(74, 14).Value = -5547.363499999999  # N74

# Row 77
$ws.Cells.Item(77, 8).Value = 5081.409  # H77
$ws.Cells.Item(77, 9).Value = 6363.4546  # I77
$ws.Cells.Item(77, 10).Value = 3799.3635  # J77
$ws.Cells.Item(77, 11).Value = 31817.273  # K77
$ws.Cells.Item(77, 12).Value = 18996.8175  # L77
$ws.Cells.Item(77, 13).Value = -27449.273  # M77
$ws.Cells.Item(77, 14).Value = -27732.8175  # N77

# Row 110
$ws.Cells.Item(110, 8).Value = 5623.067  # H110
$ws.Cells.Item(110, 9).Value = 2768.375  # I110
$ws.Cells.Item(110, 10).Value = 6661.136  # J110
$ws.Cells.Item(110, 11).Value = 2768.375  # K110
$ws.Cells.Item(110, 12).Value = 6661.136  # L110
$ws.Cells.Item(110, 13).Value = -723.375  # M110
$ws.Cells.Item(110, 14).Value = -10751.136  # N110

# Row 122
$ws.Cells.Item(122, 8).Value = 3411.4119  # H122
$ws.Cells.Item(122, 9).Value = 3228.6667  # I122
$ws.Cells.Item(122, 10).Value = 3850  # J122
$ws.Cells.Item(122, 11).Value = 9686.000100000001  # K122
$ws.Cells.Item(122, 12).Value = 11550  # L122
$ws.Cells.Item(122, 13).Value = -7236.000100000001  # M122
$ws.Cells.Item(122, 14).Value = -16450  # N122

# Row 132
$ws.Cells.Item(132, 8).Value = 583032.2  # H132
$ws.Cells.Item(132, 9).Value = 607407.4  # I132
$ws.Cells.Item(132, 10).Value = 235686.25  # J132
$ws.Cells.Item(132, 11).Value = 1822222.2  # K132
$ws.Cells.Item(132, 12).Value = 707058.75  # L132
$ws.Cells.Item(132, 13).Value = -1819692.2  # M132
$ws.Cells.Item(132, 14).Value = -712118.75  # N132


$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Cells.Item(22, 8).Value = 1713.1538  # H22
$ws.Cells.Item(22, 9).Value = 176.45454  # I22
$ws.Cells.Item(22, 10).Value = 10165  # J22
$ws.Cells.Item(22, 11).Value = 176.45454  # K22
$ws.Cells.Item(22, 12).Value = 10165  # L22
$ws.Cells.Item(22, 13).Value = -3.454540000000009  # M22
$ws.Cells.Item(22, 14).Value = -10511  # N22

# Row 105
$ws.Cells.Item(105, 8).Value = 4661.6  # H105
$ws.Cells.Item(105, 9).Value = 4769.5  # I105
$ws.Cells.Item(105, 10).Value = 4230  # J105
$ws.Cells.Item(105, 11).Value = 4769.5  # K105
$ws.Cells.Item(105, 12).Value = 4230  # L105
$ws.Cells.Item(105, 13).Value = -3022.5  # M105
$ws.Cells.Item(105, 14).Value = -7724  # N105

# Row 134
$ws.Cells.Item(134, 8).Value = 761522.8  # H134
$ws.Cells.Item(134, 9).Value = 880533.9399999999  # I134
$ws.Cells.Item(134, 10).Value = 7785.6665  # J134
$ws.Cells.Item(134, 11).Value = 2641601.82  # K134
$ws.Cells.Item(134, 12).Value = 23356.9995  # L134
$ws.Cells.Item(134, 13).Value = -2639066.82  # M134
$ws.Cells.Item(134, 14).Value = -28426.9995  # N134


$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 8659.625  # H7
$ws.Cells.Item(7, 9).Value = 11447.833  # I7
$ws.Cells.Item(7, 10).Value = 295  # J7
$ws.Cells.Item(7, 11).Value = 11447.833  # K7
$ws.Cells.Item(7, 12).Value = 295  # L7
$ws.Cells.Item(7, 13).Value = -11334.833  # M7
$ws.Cells.Item(7, 14).Value = -521  # N7

# Row 31
$ws.Cells.Item(31, 8).Value = 3884.5  # H31
$ws.Cells.Item(31, 9).Value = 1155.1875  # I31
$ws.Cells.Item(31, 10).Value = 8251.4  # J31
$ws.Cells.Item(31, 11).Value = 1155.1875  # K31
$ws.Cells.Item(31, 12).Value = 8251.4  # L31
$ws.Cells.Item(31, 13).Value = -860.1875  # M31
$ws.Cells.Item(31, 14).Value = -8841.4  # N31

# Row 34
$ws.Cells.Item(34, 8).Value = 3884.5  # H34
$ws.Cells.Item(34, 9).Value = 1155.1875  # I34
$ws.Cells.Item(34, 10).Value = 8251.4  # J34
$ws.Cells.Item(34, 11).Value = 1155.1875  # K34
$ws.Cells.Item(34, 12).Value = 8251.4  # L34
$ws.Cells.Item(34, 13).Value = -953.1875  # M34
$ws.Cells.Item(34, 14).Value = -8655.4  # N34

# Row 132
$ws.Cells.Item(132, 8).Value = 21429.066  # H132
$ws.Cells.Item(132, 9).Value = 9662.5  # I132
$ws.Cells.Item(132, 10).Value = 29273.445  # J132
$ws.Cells.Item(132, 11).Value = 28987.5  # K132
$ws.Cells.Item(132, 12).Value = 87820.33499999999  # L132
$ws.Cells.Item(132, 13).Value = -26457.5  # M132
$ws.Cells.Item(132, 14).Value = -92880.33499999999  # N132


$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Cells.Item(109, 8).Value = 10838.728  # H109
$ws.Cells.Item(109, 9).Value = 4871  # I109
$ws.Cells.Item(109, 10).Value = 18000  # J109
$ws.Cells.Item(109, 11).Value = 14613  # K109
$ws.Cells.Item(109, 12).Value = 54000  # L109
$ws.Cells.Item(109, 13).Value = -13573  # M109
$ws.Cells.Item(109, 14).Value = -56080  # N109

# Row 117
$ws.Cells.Item(117, 8).Value = 2843.3333  # H117
$ws.Cells.Item(117, 9).Value = 779  # I117
$ws.Cells.Item(117, 10).Value = 3433.1428  # J117
$ws.Cells.Item(117, 11).Value = 2337  # K117
$ws.Cells.Item(117, 12).Value = 10299.4284  # L117
$ws.Cells.Item(117, 13).Value = 1105  # M117
$ws.Cells.Item(117, 14).Value = -17183.4284  # N117

# Row 122
$ws.Cells.Item(122, 8).Value = 5524.905  # H122
$ws.Cells.Item(122, 9).Value = 1158  # I122
$ws.Cells.Item(122, 10).Value = 6252.722  # J122
$ws.Cells.Item(122, 11).Value = 10422  # K122
$ws.Cells.Item(122, 12).Value = 56274.498  # L122
$ws.Cells.Item(122, 13).Value = -7972  # M122
$ws.Cells.Item(122, 14).Value = -61174.498  # N122


$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 6155.0527  # H70
$ws.Cells.Item(70, 9).Value = 6085.727  # I70
$ws.Cells.Item(70, 10).Value = 6250.375  # J70
$ws.Cells.Item(70, 11).Value = 6085.727  # K70
$ws.Cells.Item(70, 12).Value = 6250.375  # L70
$ws.Cells.Item(70, 13).Value = -5815.727  # M70
$ws.Cells.Item(70, 14).Value = -6790.375  # N70

# Row 73
$ws.Cells.Item(73, 8).Value = 6155.0527  # H73
$ws.Cells.Item(73, 9).Value = 6085.727  # I73
$ws.Cells.Item(73, 10).Value = 6250.375  # J73
$ws.Cells.Item(73, 11).Value = 6085.727  # K73
$ws.Cells.Item(73, 12).Value = 6250.375  # L73
$ws.Cells.Item(73, 13).Value = -5149.727  # M73
$ws.Cells.Item(73, 14).Value = -8122.375  # N73

# Row 122
$ws.Cells.Item(122, 8).Value = 5688.8125  # H122
$ws.Cells.Item(122, 9).Value = 4052.9092  # I122
$ws.Cells.Item(122, 10).Value = 9287.799999999999  # J122
$ws.Cells.Item(122, 11).Value = 12158.7276  # K122
$ws.Cells.Item(122, 12).Value = 27863.4  # L122
$ws.Cells.Item(122, 13).Value = -9708.7276  # M122
$ws.Cells.Item(122, 14).Value = -32763.4  # N122

# Row 132
$ws.Cells.Item(132, 8).Value = 8443.1  # H132
$ws.Cells.Item(132, 9).Value = 9571.200000000001  # I132
$ws.Cells.Item(132, 10).Value = 5058.8  # J132
$ws.Cells.Item(132, 11).Value = 28713.6  # K132
$ws.Cells.Item(132, 12).Value = 15176.4  # L132
$ws.Cells.Item(132, 13).Value = -26183.6  # M132
$ws.Cells.Item(132, 14).Value = -20236.4  # N132


$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Cells.Item(61, 8).Value = 6388.36  # H61
$ws.Cells.Item(61, 9).Value = 5293.1055  # I61
$ws.Cells.Item(61, 10).Value = 9856.666999999999  # J61
$ws.Cells.Item(61, 11).Value = 5293.1055  # K61
$ws.Cells.Item(61, 12).Value = 9856.666999999999  # L61
$ws.Cells.Item(61, 13).Value = -5091.1055  # M61
$ws.Cells.Item(61, 14).Value = -10260.667  # N61

# Row 100
$ws.Cells.Item(100, 8).Value = 3024.7307  # H100
$ws.Cells.Item(100, 9).Value = 3905.4546  # I100
$ws.Cells.Item(100, 10).Value = 2378.8667  # J100
$ws.Cells.Item(100, 11).Value = 3905.4546  # K100
$ws.Cells.Item(100, 12).Value = 2378.8667  # L100
$ws.Cells.Item(100, 13).Value = -3364.4546  # M100
$ws.Cells.Item(100, 14).Value = -3460.8667  # N100

# Row 113
$ws.Cells.Item(113, 8).Value = 6388.36  # H113
$ws.Cells.Item(113, 9).Value = 5293.1055  # I113
$ws.Cells.Item(113, 10).Value = 9856.666999999999  # J113
$ws.Cells.Item(113, 11).Value = 5293.1055  # K113
$ws.Cells.Item(113, 12).Value = 9856.666999999999  # L113
$ws.Cells.Item(113, 13).Value = -3123.1055  # M113
$ws.Cells.Item(113, 14).Value = -14196.667  # N113

# Row 132
$ws.Cells.Item(132, 8).Value = 1951  # H132
$ws.Cells.Item(132, 9).Value = 1951  # I132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 11).Value = 5853  # K132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 13).Value = -3323  # M132
$ws.Cells.Item(132, 14).ClearContents()  # N132

# Row 136
$ws.Cells.Item(136, 8).Value = 41675520  # H136
$ws.Cells.Item(136, 9).Value = 83342344  # I136
$ws.Cells.Item(136, 10).Value = 8696.416999999999  # J136
$ws.Cells.Item(136, 11).Value = 250027032  # K136
$ws.Cells.Item(136, 12).Value = 26089.251  # L136
$ws.Cells.Item(136, 13).Value = -250024482  # M136
$ws.Cells.Item(136, 14).Value = -31189.251  # N136

# Row 139
$ws.Cells.Item(139, 8).Value = 149000  # H139
$ws.Cells.Item(139, 9).Value = 0  # I139
$ws.Cells.Item(139, 10).Value = 149000  # J139
$ws.Cells.Item(139, 11).Value = 0  # K139
$ws.Cells.Item(139, 12).Value = 149000  # L139
$ws.Cells.Item(139, 14).Value = -159280  # N139


$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 5526.773  # H62
$ws.Cells.Item(62, 9).Value = 7970.4  # I62
$ws.Cells.Item(62, 10).Value = 4808.0586  # J62
$ws.Cells.Item(62, 11).Value = 7970.4  # K62
$ws.Cells.Item(62, 12).Value = 4808.0586  # L62
$ws.Cells.Item(62, 13).Value = -7346.4  # M62
$ws.Cells.Item(62, 14).Value = -6056.0586  # N62

# Row 65
$ws.Cells.Item(65, 8).Value = 5526.773  # H65
$ws.Cells.Item(65, 9).Value = 7970.4  # I65
$ws.Cells.Item(65, 10).Value = 4808.0586  # J65
$ws.Cells.Item(65, 11).Value = 39852  # K65
$ws.Cells.Item(65, 12).Value = 24040.293  # L65
$ws.Cells.Item(65, 13).Value = -36732  # M65
$ws.Cells.Item(65, 14).Value = -30280.293  # N65

# Row 107
$ws.Cells.Item(107, 8).Value = 4598.6  # H107
$ws.Cells.Item(107, 9).Value = 5166.6665  # I107
$ws.Cells.Item(107, 10).Value = 3746.5  # J107
$ws.Cells.Item(107, 11).Value = 15499.9995  # K107
$ws.Cells.Item(107, 12).Value = 11239.5  # L107
$ws.Cells.Item(107, 13).Value = -13579.9995  # M107
$ws.Cells.Item(107, 14).Value = -15079.5  # N107

# Row 122
$ws.Cells.Item(122, 8).Value = 17804.445  # H122
$ws.Cells.Item(122, 9).Value = 11290.583  # I122
$ws.Cells.Item(122, 10).Value = 30832.166  # J122
$ws.Cells.Item(122, 11).Value = 33871.749  # K122
$ws.Cells.Item(122, 12).Value = 92496.49800000001  # L122
$ws.Cells.Item(122, 13).Value = -31421.749  # M122
$ws.Cells.Item(122, 14).Value = -97396.49800000001  # N122

# Row 132
$ws.Cells.Item(132, 8).Value = 9584.526  # H132
$ws.Cells.Item(132, 9).Value = 5867.7407  # I132
$ws.Cells.Item(132, 10).Value = 18707.545  # J132
$ws.Cells.Item(132, 11).Value = 17603.2221  # K132
$ws.Cells.Item(132, 12).Value = 56122.63499999999  # L132
$ws.Cells.Item(132, 13).Value = -15073.2221  # M132
$ws.Cells.Item(132, 14).Value = -61182.63499999999  # N132

# Row 136
$ws.Cells.Item(136, 8).Value = 13167491  # H136
$ws.Cells.Item(136, 9).Value = 16136978  # I136
$ws.Cells.Item(136, 10).Value = 16908  # J136
$ws.Cells.Item(136, 11).Value = 48410934  # K136
$ws.Cells.Item(136, 12).Value = 50724  # L136
$ws.Cells.Item(136, 13).Value = -48408384  # M136
$ws.Cells.Item(136, 14).Value = -55824  # N136
